$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(according to the population census data)" note row (old row 2)
$ws.Rows.Item(2).Delete()

# Drop the 1989 and 2002 year columns, keeping only the 2014 figure
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# Match the updated row heights used in the trimmed layout
$ws.Rows.Item(1).RowHeight = 20.1
$ws.Rows.Item(2).RowHeight = 20.1
$ws.Rows.Item(3).RowHeight = 20.1
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
